$d = $word.ActiveDocument

# The "Final Decision" section (heading + Input/Output lines) is being
# removed entirely; the blank paragraph that precedes it (and the document's
# trailing blank paragraph) stay untouched.

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*6. *Final Decision*") {
        $startPara = $p
    }
    if ($t -like "*Output: Personalized story ending*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
